$wb = $excel.ActiveWorkbook

$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

$status = "Handed back: in sync with en-US"
$zhTargetFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$deTargetFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$zhHandbackTime = "2016-08-20 06:44:10"
$deHandbackTime = "2016-08-20 06:44:16"
$aLink = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1363d3dd7dc2ba43ee3f7eba0a67c16ffba220ac/e2e/a.md"
$bLink = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1363d3dd7dc2ba43ee3f7eba0a67c16ffba220ac/e2e/b.md"

# Overview sheet status columns (E = zh-cn, F = de-de)
$overview.Range("E2").Value = $status
$overview.Range("F2").Value = $status
$overview.Range("E3").Value = $status
$overview.Range("F3").Value = $status

# zh-cn sheet updates
$zh.Range("C2").Value = $status
$zh.Range("C3").Value = $status
$zh.Range("I2").Value = "a.md"
$zh.Range("I3").Value = "a.md"
$zh.Range("J2").Value = $zhTargetFile
$zh.Range("J3").Value = $zhTargetFile
$zh.Range("K2").Value = $zhHandbackTime
$zh.Range("K3").Value = $zhHandbackTime

$zh.Hyperlinks.Add($zh.Range("I2"), $aLink, "", "", "a.md")
$zh.Hyperlinks.Add($zh.Range("I3"), $aLink, "", "", "a.md")

# de-de sheet updates
$de.Range("C2").Value = $status
$de.Range("C3").Value = $status
$de.Range("I2").Value = "a.md"
$de.Range("I3").Value = "a.md"
$de.Range("J2").Value = $deTargetFile
$de.Range("J3").Value = $deTargetFile
$de.Range("K2").Value = $deHandbackTime
$de.Range("K3").Value = $deHandbackTime

$de.Hyperlinks.Add($de.Range("I2"), $aLink, "", "", "a.md")
$de.Hyperlinks.Add($de.Range("I3"), $aLink, "", "", "a.md")
